# Commit: "Fixed update to excel issue"
# - Rename header B1 in "Weekly Quantity" and "Monthly Trend" sheets
# - Add a new "PO Forecast" sheet with ds / PO_Forecast / yhat_lower / yhat_upper data

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

$weekly.Range("B1").Value = "Weekly_PO_Qty"
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Copy header formatting (bold + border + alignment) from an existing header row
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)

# Copy date-column formatting (custom date/time number format) from an existing date column
$weekly.Range("A2").Copy()
$forecast.Range("A2:A46").PasteSpecial(-4122)

$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(45165.99999999999, 23, -18.52974402814541, 64.05857858370943),
    @(45172.99999999999, 23, -17.0624433464677, 61.51253670614904),
    @(45179.99999999999, 23, -17.29829642463121, 66.24270436454208),
    @(45186.99999999999, 23, -18.38585016383175, 62.31256070168083),
    @(45193.99999999999, 23, -16.86899270789769, 66.7599329940746),
    @(45277.99999999999, 23, -16.85410279866685, 66.55464358019586),
    @(45298.99999999999, 23, -18.7825430914603, 65.75632400690988),
    @(45319.99999999999, 24, -18.26646437591452, 63.98425490549199),
    @(45326.99999999999, 24, -17.06335660634159, 63.5276773473262),
    @(45333.99999999999, 24, -19.6495681130787, 64.30944647050447),
    @(45340.99999999999, 24, -16.90354989162013, 65.15783471109923),
    @(45347.99999999999, 24, -16.4725983532117, 65.42990749343181),
    @(45354.99999999999, 24, -15.70251757121314, 65.74697868901647),
    @(45361.99999999999, 24, -12.18692881790845, 64.90252912732456),
    @(45368.99999999999, 24, -20.78293438779249, 63.86820825548789),
    @(45375.99999999999, 24, -18.59802344109837, 62.05104717237122),
    @(45410.99999999999, 24, -16.849055511509, 66.5457943700234),
    @(45417.99999999999, 24, -17.22037182636603, 66.40066505643563),
    @(45424.99999999999, 24, -16.20172195055751, 67.18324854661735),
    @(45438.99999999999, 24, -18.03679687244118, 62.91136057935725),
    @(45445.99999999999, 24, -15.84012486811555, 61.99013993032371),
    @(45515.99999999999, 25, -15.05817769481522, 65.0105218699664),
    @(45522.99999999999, 25, -14.19225276290955, 64.29088629198861),
    @(45536.99999999999, 25, -15.24433244442476, 65.48390213461551),
    @(45543.99999999999, 25, -14.83883588828091, 60.12598780554303),
    @(45557.99999999999, 25, -16.32181881479075, 66.90932210573321),
    @(45564.99999999999, 25, -19.30969250477439, 64.87623654467902),
    @(45578.99999999999, 25, -14.23712313274605, 66.83316783239889),
    @(45585.99999999999, 25, -16.04739533331695, 65.31917125126057),
    @(45592.99999999999, 25, -15.07245908288202, 64.94708621970206),
    @(45599.99999999999, 25, -16.05483135034347, 64.48507928765889),
    @(45606.99999999999, 25, -14.970657724795, 68.90952998207514),
    @(45613.99999999999, 25, -13.81684398988979, 65.15711520070462),
    @(45620.99999999999, 25, -14.45742597512262, 66.27765068853255),
    @(45627.99999999999, 25, -16.47184180299771, 62.65227906367382),
    @(45634.99999999999, 25, -17.45730058343292, 63.45421303427214),
    @(45641.99999999999, 25, -14.64445523199103, 65.7076484358786),
    @(45648.99999999999, 25, -16.03203287054095, 67.12440826009228),
    @(45655.99999999999, 25, -19.38771768242889, 66.27624628844353),
    @(45662.99999999999, 25, -13.4183635145329, 62.35471349545814),
    @(45669.99999999999, 25, -12.29905766149143, 63.92959044933119),
    @(45676.99999999999, 25, -14.64919378538262, 66.33912232621901),
    @(45683.99999999999, 25, -13.47659214107084, 63.39000003963893),
    @(45690.99999999999, 25, -14.16153437036525, 66.28540286811558),
    @(45697.99999999999, 26, -15.75191973107698, 67.77106326429301)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $forecast.Cells.Item($r, 1).Value = $row[0]
    $forecast.Cells.Item($r, 2).Value = $row[1]
    $forecast.Cells.Item($r, 3).Value = $row[2]
    $forecast.Cells.Item($r, 4).Value = $row[3]
}
